$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J13").Value = "5:30 am meet for Carlie K/Trevor M at Grafton Park n Ride"
$ws.Range("F14").Value = "SMALL CHANGES-RESET`n4:45 am meet for Carlie K/Trevor M at Grafton Park n Ride"
$ws.Range("N14").Value = "SET UP ON REG FAR RIGHT BY OFFICE`n5:15 am meet for Carlie K/Trevor M at Grafton Park n Ride`t"
$ws.Range("R14").Value = "4:45 am meet for Carlie K/Trevor M at Grafton Park n Ride"
$ws.Range("J15").Value = "Brianna H"
$ws.Range("V15").Value = "DJ S"
$ws.Range("F16").Value = "Kim G"
$ws.Range("J16").Value = "Aaron M"
$ws.Range("N16").Value = "Nate C"
$ws.Range("R16").Value = "DJ S"
$ws.Range("V16").Value = "Aaron M"
$ws.Range("F17").Value = "Carlie K"
$ws.Range("G17").Value = "@ Store`n(w/ Trevor M)"
$ws.Range("J17").Value = "Carlie K"
$ws.Range("K17").Value = "@ Store`n(w/ Trevor M)"
$ws.Range("N17").Value = "Aaron M"
$ws.Range("R17").Value = "Carlie K"
$ws.Range("S17").Value = "@ Store`n(w/ Trevor M)"
$ws.Range("V17").Value = "Greg H"
$ws.Range("F18").Value = "Curt B"
$ws.Range("J18").Value = "Greg H"
$ws.Range("N18").Value = "Carlie K"
$ws.Range("O18").Value = "@ Store`n(w/ Trevor M)"
$ws.Range("R18").Value = "Monica G"
$ws.Range("F19").Value = "Cynthia M"
$ws.Range("J19").Value = "Monica G"
$ws.Range("N19").Value = "Cynthia M"
$ws.Range("R19").Value = "Stephanie G"
$ws.Range("V19").Value = "Sonia T"
$ws.Range("F20").Value = "Greg H"
$ws.Range("J20").Value = "Robyn K"
$ws.Range("N20").Value = "Greg H"
$ws.Range("R20").Value = "Trevor M"
$ws.Range("S20").Value = "@ Store`n(w/ Carlie K)"
$ws.Range("F21").Value = "Robyn K"
$ws.Range("J21").Value = "Sonia T"
$ws.Range("N21").Value = "Mai M"
$ws.Range("V21").Value = "Evelin A"
$ws.Range("F22").Value = "Sonia T"
$ws.Range("J22").Value = "Stephanie G"
$ws.Range("N22").Value = "Monica G"
$ws.Range("V22").Value = "Qiana B"
$ws.Range("F23").Value = "Trevor M"
$ws.Range("G23").Value = "@ Store`n(w/ Carlie K)"
$ws.Range("J23").Value = "Sue M"
$ws.Range("N23").Value = "Paul T"
$ws.Range("J24").Value = "Trevor M"
$ws.Range("K24").Value = "@ Store`n(w/ Carlie K)"
$ws.Range("N24").Value = "Sonia T"
$ws.Range("N25").Value = "Stephanie G"
$ws.Range("N26").Value = "Trevor M"
$ws.Range("O26").Value = "@ Store`n(w/ Carlie K)"
$ws.Range("N27").Value = "Angela B"
$ws.Range("N28").Value = "Anisha V"
$ws.Range("N29").Value = "Brianna H"
$ws.Range("N30").Value = "Evelin A"
$ws.Range("R30").Value = "Brianna H"
$ws.Range("N31").Value = "Joseph S"
$ws.Range("R31").Value = "Casey V"
$ws.Range("N32").Value = "Justin L"
$ws.Range("R32").Value = "Lashaun C"
$ws.Range("V32").Value = "Nate C"
$ws.Range("R33").Value = "Lori N"
$ws.Range("V33").Value = "Kim G"
$ws.Range("N34").Value = "Qiana B"
$ws.Range("V34").Value = "Lori N"
$ws.Range("J35").Value = "Kim G"
$ws.Range("N35").Value = "Taylor G"
$ws.Range("J36").Value = "Cynthia M"
$ws.Range("J37").Value = "Mai M"
$ws.Range("F38").Value = "Ian K"
$ws.Range("F39").Value = "Paul T"
$ws.Range("R43").Value = "Kim G"
$ws.Range("R44").Value = "Cynthia M"
$ws.Range("V44").Value = "Lashaun C"
$ws.Range("R45").Value = "Paul T"
$ws.Range("V45").Value = "Carlie K"
$ws.Range("R46").Value = "Sonia T"
$ws.Range("V46").Value = "Casey V"
$ws.Range("V47").Value = "Stephanie G"
$ws.Range("J52").Value = "Ian K"
$ws.Range("J53").Value = "Paul T"
$ws.Range("N53").Value = "Katherine G"
$ws.Range("N54").Value = "Ian K"
$ws.Range("V56").Value = "Katherine G"
$ws.Range("V57").Value = "Cynthia M"
$ws.Range("V58").Value = "Mai M"
$ws.Range("F59").Value = "Katherine G"
$ws.Range("V59").Value = "Monica G"
$ws.Range("F60").Value = "Sue M"
$ws.Range("R61").Value = "Katherine G"
$ws.Range("R62").Value = "Ian K"
$ws.Range("V65").Value = "ROTE OIL #14 Trevor M (CITGO)"
$ws.Range("V70").Value = "ROTE OIL #13 Trevor M (BP)"
$ws.Range("R72").Value = "Nate C"
$ws.Range("R73").Value = "Curt B"
$ws.Range("N74").Value = "Sarah H"
$ws.Range("V74").Value = "Ian K"
$ws.Range("F75").Value = "DJ S"
$ws.Range("R75").Value = "Robyn K"
$ws.Range("V75").Value = "Paul T"
$ws.Range("F76").Value = "Mai M"
$ws.Range("N76").Value = "Eva G"
$ws.Range("R76").Value = "Sue M"
$ws.Range("J77").Value = "Katherine G"
$ws.Range("N77").Value = "Lori N"
$ws.Range("R77").Value = "Evelin A"
$ws.Range("J78").Value = "Casey V"
$ws.Range("R78").Value = "Jeri H"
$ws.Range("F79").Value = "IL: 4:40 AM MEET AT JANESVILLE PARK N RIDE"
$ws.Range("J79").Value = "Lashaun C"
$ws.Range("R79").Value = "Justin L"
$ws.Range("N81").Value = "Kim G"
$ws.Range("R81").Value = "Taylor G"
$ws.Range("N82").Value = "Lashaun C"
$ws.Range("J89").Value = "Nate C"
$ws.Range("V89").Value = "Sarah H"
$ws.Range("J90").Value = "Angela B"
$ws.Range("V90").Value = "Brianna H"
$ws.Range("J91").Value = "Anisha V"
$ws.Range("R91").Value = "Sarah H"
$ws.Range("F92").Value = "Brianna H"
$ws.Range("J92").Value = "Elijah E"
$ws.Range("R92").Value = "Angela B"
$ws.Range("F93").Value = "Josie N"
$ws.Range("J93").Value = "Evelin A"
$ws.Range("R93").Value = "Anisha V"
$ws.Range("F94").Value = "Lori N"
$ws.Range("J94").Value = "Joseph S"
$ws.Range("F95").Value = "Qiana B"
$ws.Range("J95").Value = "Taylor G"
$ws.Range("R95").Value = "Elijah E"
$ws.Range("R96").Value = "Eva G"
$ws.Range("R97").Value = "Joseph S"
$ws.Range("R99").Value = "Qiana B"
$ws.Range("F105").Value = "Sarah H"
$ws.Range("F106").Value = "Angela B"
$ws.Range("F107").Value = "Anisha V"
$ws.Range("F109").Value = "Elijah E"
$ws.Range("F110").Value = "Eva G"
$ws.Range("J110").Value = "Sarah H"
$ws.Range("F111").Value = "Evelin A"
$ws.Range("F112").Value = "Joseph S"
$ws.Range("J112").Value = "Eva G"
$ws.Range("F113").Value = "Justin L"
$ws.Range("J113").Value = "Josie N"
$ws.Range("J114").Value = "Lori N"
$ws.Range("F115").Value = "Nate C"
$ws.Range("F116").Value = "Taylor G"
